# The workbook has a single sheet ("Avverkningsanmälningar").
# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2023-09-19 -> 2023-09-20, i.e. 45188 -> 45189) for every data row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value()
    if ($current -ne $null) {
        $cell.Value = $current.AddDays(1)
    }
}
